# Generate Report for Handoff
# Refresh the "Latest Handoff Date"/"Latest Handoff Datetime" for the file
# 36e46dc5-c445-48ce-af7f-65d1a668a69e.md which has just been handed off again.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-32-19 16:32:27"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-19 16:32:25"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-19 16:32:27"
